# Add season-record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column headers, styled like the existing header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-41: season record values for each player row.
for ($row = 2; $row -le 41; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 0
}
